# Mise à jour de l'auto-évaluation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New remark cells in column C ---
# Row 16: add remark "A faire par Vanessa" (default style)
$ws.Range("C16").Value = "A faire par Vanessa"

# Row 18 & 20: add remark "Presque bon" (default style)
$ws.Range("C18").Value = "Presque bon"
$ws.Range("C20").Value = "Presque bon"

# Row 32 & 33: add remark "A FAIRE !!" with bold style (same cellXf as D1/B26, s="3")
$ws.Range("C32").Value = "A FAIRE !!"
$ws.Range("C33").Value = "A FAIRE !!"
$ws.Range("D1").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("C33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 34: remove the remark "Vérifier que c'est correct" from C34 ---
$ws.Range("C34").ClearContents()

# --- Fill-color (style) changes on column B ---
# B17 and B19: change fill from style 10 (orange) to style 9 (theme accent), copy format from B12
$ws.Range("B12").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)

# B34 and B35: change fill from style 11 (yellow) to style 9 (theme accent), copy format from B12
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("B35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Sheet view: update scroll position and active cell selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("B41").Select()
